$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new notes for rows 13, 14, 24, 27, 44, 45 in column B
$ws.Range("B13").Value = "tcm0001, tcm0002 - The test doesn't really check the function, only that some output is read"
$ws.Range("B14").Value = "tcm0001, tcm0002 - The test doesn't really check the function, only that we get something"
$ws.Range("B24").Value = "Has only been tested by running it manually, not tested by test cases"
$ws.Range("B27").Value = "tcm0001"
$ws.Range("B45").Value = "tcm0001 - does not test all aspects of the adapter, but some."
$ws.Range("B44").Value = "tcm0002 - does not test all aspects of the adapter, but some."

# Update the view: scroll so row 25 is the top-left visible row, and select B43
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("B43").Select()
